$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "245.45"
Set-TextValue "E2" "-0.70%"
Set-TextValue "D3" "27.20"
Set-TextValue "E3" "2.94%"
Set-TextValue "D4" "5.101"
Set-TextValue "E4" "0.68%"
Set-TextValue "D5" "0.05703"
Set-TextValue "E5" "1.74%"
Set-TextValue "E6" "0.49%"
Set-TextValue "D7" "0.8190"
Set-TextValue "E7" "0.72%"
Set-TextValue "D8" "0.8567"
Set-TextValue "E8" "1.58%"
Set-TextValue "D9" "0.06944"
Set-TextValue "E9" "-0.81%"
Set-TextValue "D10" "0.02839"
Set-TextValue "E10" "-0.59%"
Set-TextValue "D11" "0.09393"
Set-TextValue "E11" "-0.09%"
Set-TextValue "D12" "0.001517"
Set-TextValue "E12" "0.14%"
Set-TextValue "D13" "0.04044"
Set-TextValue "E13" "-13.08%"
Set-TextValue "D14" "0.0006005"
Set-TextValue "E14" "0.77%"
Set-TextValue "D15" "0.006210"
Set-TextValue "E15" "0.74%"
Set-TextValue "E16" "-2.65%"
Set-TextValue "D17" "3.007"
Set-TextValue "E17" "-0.22%"
Set-TextValue "D18" "2.229"
Set-TextValue "E18" "8.44%"
Set-TextValue "D20" "0.1332"
Set-TextValue "E20" "-0.49%"
Set-TextValue "D21" "0.03221"
Set-TextValue "E21" "0.69%"
Set-TextValue "E22" "-1.83%"
Set-TextValue "D23" "3.593"
Set-TextValue "E23" "-3.89%"
Set-TextValue "D25" "0.001218"
Set-TextValue "E25" "-2.03%"
Set-TextValue "D26" "0.004474"
Set-TextValue "E26" "-2.47%"
Set-TextValue "D27" "0.00009896"
Set-TextValue "E28" "3.59%"
Set-TextValue "D40" "0.03732"
Set-TextValue "E40" "1.78%"
Set-TextValue "D41" "0.005988"
Set-TextValue "E41" "-2.64%"
Set-TextValue "D42" "0.1060"
Set-TextValue "E42" "0.20%"
Set-TextValue "D43" "0.002439"
Set-TextValue "E43" "-2.43%"
Set-TextValue "D44" "0.009689"
Set-TextValue "E44" "17.28%"
Set-TextValue "D45" "0.00005152"
Set-TextValue "E45" "-4.55%"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "E46" "-0.08%"
Set-TextValue "E47" "-8.25%"
Set-TextValue "E48" "-3.66%"
Set-TextValue "D49" "0.00002098"
Set-TextValue "E49" "-0.08%"
Set-TextValue "D50" "0.0001998"
Set-TextValue "E50" "-0.08%"

Write-Host "Updated 65 cells with new symbol list data"
